$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 316.20834
$ws.Range("I28").Value = 304.45
$ws.Range("J28").Value = 375
$ws.Range("K28").Value = 304.45
$ws.Range("L28").Value = 375
$ws.Range("M28").Value = 180.55
$ws.Range("N28").Value = -1345
$ws.Range("H58").Value = 1061.8462
$ws.Range("I58").Value = 454
$ws.Range("J58").Value = 1582.8572
$ws.Range("K58").Value = 1362
$ws.Range("L58").Value = 4748.571599999999
$ws.Range("M58").Value = -1212
$ws.Range("N58").Value = -5048.571599999999
$ws.Range("H62").Value = 2870.0435
$ws.Range("I62").Value = 4913.125
$ws.Range("J62").Value = 1780.4
$ws.Range("K62").Value = 4913.125
$ws.Range("L62").Value = 1780.4
$ws.Range("M62").Value = -4289.125
$ws.Range("N62").Value = -3028.4
$ws.Range("H65").Value = 2870.0435
$ws.Range("I65").Value = 4913.125
$ws.Range("J65").Value = 1780.4
$ws.Range("K65").Value = 24565.625
$ws.Range("L65").Value = 8902
$ws.Range("M65").Value = -21445.625
$ws.Range("N65").Value = -15142
$ws.Range("H75").Value = 43345.715
$ws.Range("J75").Value = 43345.715
$ws.Range("L75").Value = 43345.715
$ws.Range("N75").Value = -45217.715
$ws.Range("H76").Value = 3237.182
$ws.Range("I76").Value = 3002.5
$ws.Range("J76").Value = 3289.3333
$ws.Range("K76").Value = 3002.5
$ws.Range("L76").Value = 3289.3333
$ws.Range("M76").Value = -2687.5
$ws.Range("N76").Value = -3919.3333
$ws.Range("H78").Value = 43345.715
$ws.Range("J78").Value = 43345.715
$ws.Range("L78").Value = 130037.145
$ws.Range("N78").Value = -139397.145
$ws.Range("H79").Value = 3237.182
$ws.Range("I79").Value = 3002.5
$ws.Range("J79").Value = 3289.3333
$ws.Range("K79").Value = 3002.5
$ws.Range("L79").Value = 3289.3333
$ws.Range("M79").Value = -1910.5
$ws.Range("N79").Value = -5473.3333
$ws.Range("H86").Value = 12512618
$ws.Range("I86").Value = 25005260
$ws.Range("J86").Value = 19976
$ws.Range("K86").Value = 25005260
$ws.Range("L86").Value = 19976
$ws.Range("M86").Value = -25004137
$ws.Range("N86").Value = -22222
$ws.Range("H89").Value = 12512618
$ws.Range("I89").Value = 25005260
$ws.Range("J89").Value = 19976
$ws.Range("K89").Value = 125026300
$ws.Range("L89").Value = 99880
$ws.Range("M89").Value = -125020684
$ws.Range("N89").Value = -111112
$ws.Range("H92").Value = 1438.375
$ws.Range("I92").Value = 1329.7142
$ws.Range("J92").Value = 2199
$ws.Range("K92").Value = 1329.7142
$ws.Range("L92").Value = 2199
$ws.Range("M92").Value = -81.71419999999989
$ws.Range("N92").Value = -4695
$ws.Range("H98").Value = 32532.8
$ws.Range("I98").Value = 924.2308
$ws.Range("J98").Value = 66775.414
$ws.Range("K98").Value = 924.2308
$ws.Range("L98").Value = 66775.414
$ws.Range("M98").Value = 573.7692
$ws.Range("N98").Value = -69771.414
$ws.Range("H106").Value = 112468
$ws.Range("I106").Value = 1399.8
$ws.Range("J106").Value = 390138.5
$ws.Range("K106").Value = 1399.8
$ws.Range("L106").Value = 390138.5
$ws.Range("M106").Value = -768.8
$ws.Range("N106").Value = -391400.5
$ws.Range("H107").Value = 6773
$ws.Range("I107").Value = 8147.5386
$ws.Range("K107").Value = 8147.5386
$ws.Range("M107").Value = -6227.5386
$ws.Range("H109").Value = 36511
$ws.Range("J109").Value = 36511
$ws.Range("L109").Value = 36511
$ws.Range("N109").Value = -39285
$ws.Range("H114").Value = 32649.5
$ws.Range("J114").Value = 32649.5
$ws.Range("L114").Value = 32649.5
$ws.Range("N114").Value = -41327.5
$ws.Range("H122").Value = 32532.8
$ws.Range("I122").Value = 924.2308
$ws.Range("J122").Value = 66775.414
$ws.Range("K122").Value = 2772.6924
$ws.Range("L122").Value = 200326.242
$ws.Range("M122").Value = -322.6923999999999
$ws.Range("N122").Value = -205226.242
$ws.Range("H124").Value = 36907.2
$ws.Range("I124").Value = 1000
$ws.Range("J124").Value = 45884
$ws.Range("K124").Value = 1000
$ws.Range("L124").Value = 45884
$ws.Range("M124").Value = 3910
$ws.Range("N124").Value = -55704
$ws.Range("H126").Value = 36794.4
$ws.Range("J126").Value = 36794.4
$ws.Range("L126").Value = 36794.4
$ws.Range("N126").Value = -46674.4
$ws.Range("H128").Value = 40527.2
$ws.Range("J128").Value = 40527.2
$ws.Range("L128").Value = 40527.2
$ws.Range("N128").Value = -50487.2
$ws.Range("H130").Value = 43398.4
$ws.Range("J130").Value = 43398.4
$ws.Range("L130").Value = 43398.4
$ws.Range("N130").Value = -53438.4
$ws.Range("H133").Value = 49863.062
$ws.Range("J133").Value = 49863.062
$ws.Range("L133").Value = 49863.062
$ws.Range("N133").Value = -59983.062

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4714.7095
$ws.Range("I2").Value = 5370.885
$ws.Range("J2").Value = 1302.6
$ws.Range("K2").Value = 5370.885
$ws.Range("L2").Value = 1302.6
$ws.Range("M2").Value = -5257.885
$ws.Range("N2").Value = -1528.6
$ws.Range("H24").Value = 28724.5
$ws.Range("J24").Value = 28724.5
$ws.Range("L24").Value = 28724.5
$ws.Range("N24").Value = -29472.5
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H45").Value = 2582.6
$ws.Range("I45").Value = 2301.7144
$ws.Range("J45").Value = 3238
$ws.Range("K45").Value = 2301.7144
$ws.Range("L45").Value = 3238
$ws.Range("M45").Value = -1924.7144
$ws.Range("N45").Value = -3992
$ws.Range("H80").Value = 48094
$ws.Range("J80").Value = 48094
$ws.Range("L80").Value = 48094
$ws.Range("N80").Value = -50090
$ws.Range("H83").Value = 48094
$ws.Range("J83").Value = 48094
$ws.Range("L83").Value = 144282
$ws.Range("N83").Value = -154266
$ws.Range("H97").Value = 742.9583
$ws.Range("I97").Value = 529.44446
$ws.Range("J97").Value = 1383.5
$ws.Range("K97").Value = 529.44446
$ws.Range("L97").Value = 1383.5
$ws.Range("M97").Value = -33.44446000000005
$ws.Range("N97").Value = -2375.5
$ws.Range("H100").Value = 28724.5
$ws.Range("J100").Value = 28724.5
$ws.Range("L100").Value = 28724.5
$ws.Range("N100").Value = -30888.5
$ws.Range("H102").Value = 18498.334
$ws.Range("I102").Value = 1627.5
$ws.Range("K102").Value = 1627.5
$ws.Range("M102").Value = -5.5
$ws.Range("H110").Value = 996.6667
$ws.Range("I110").Value = 853.5714
$ws.Range("K110").Value = 853.5714
$ws.Range("M110").Value = 1191.4286
$ws.Range("H113").Value = 37754.2
$ws.Range("J113").Value = 37754.2
$ws.Range("L113").Value = 37754.2
$ws.Range("N113").Value = -46432.2
$ws.Range("H114").Value = 39260.4
$ws.Range("J114").Value = 39260.4
$ws.Range("L114").Value = 39260.4
$ws.Range("N114").Value = -47938.4
$ws.Range("H116").Value = 4714.7095
$ws.Range("I116").Value = 5370.885
$ws.Range("J116").Value = 1302.6
$ws.Range("K116").Value = 5370.885
$ws.Range("L116").Value = 1302.6
$ws.Range("M116").Value = -3076.885
$ws.Range("N116").Value = -5890.6
$ws.Range("H117").Value = 49297.4
$ws.Range("J117").Value = 49297.4
$ws.Range("L117").Value = 49297.4
$ws.Range("N117").Value = -58475.4
$ws.Range("H118").Value = 44998
$ws.Range("J118").Value = 44998
$ws.Range("L118").Value = 44998
$ws.Range("N118").Value = -48312
$ws.Range("H119").Value = 45090.5
$ws.Range("J119").Value = 45090.5
$ws.Range("L119").Value = 45090.5
$ws.Range("N119").Value = -54766.5
$ws.Range("H121").Value = 48255
$ws.Range("J121").Value = 48255
$ws.Range("L121").Value = 48255
$ws.Range("N121").Value = -51749
$ws.Range("H122").Value = 2202.4614
$ws.Range("I122").Value = 2202.4614
$ws.Range("K122").Value = 6607.3842
$ws.Range("M122").Value = -4157.3842
$ws.Range("H123").Value = 51429
$ws.Range("J123").Value = 51429
$ws.Range("L123").Value = 51429
$ws.Range("N123").Value = -61229
$ws.Range("H128").Value = 48137.668
$ws.Range("J128").Value = 48137.668
$ws.Range("L128").Value = 48137.668
$ws.Range("N128").Value = -58097.668
$ws.Range("H131").Value = 50563.668
$ws.Range("J131").Value = 50563.668
$ws.Range("L131").Value = 50563.668
$ws.Range("N131").Value = -60643.668
$ws.Range("H132").Value = 12822490
$ws.Range("I132").Value = 20001306
$ws.Range("K132").Value = 60003918
$ws.Range("M132").Value = -60001388

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4714.7095
$ws.Range("I3").Value = 5370.885
$ws.Range("J3").Value = 1302.6
$ws.Range("K3").Value = 5370.885
$ws.Range("L3").Value = 1302.6
$ws.Range("M3").Value = -5256.885
$ws.Range("N3").Value = -1530.6
$ws.Range("H105").Value = 3617.4285
$ws.Range("I105").Value = 3862.2
$ws.Range("K105").Value = 3862.2
$ws.Range("M105").Value = -2115.2
$ws.Range("H107").Value = 2053.6667
$ws.Range("I107").Value = 1946.0385
$ws.Range("J107").Value = 2753.25
$ws.Range("K107").Value = 1946.0385
$ws.Range("L107").Value = 2753.25
$ws.Range("M107").Value = -26.03850000000011
$ws.Range("N107").Value = -6593.25
$ws.Range("H111").Value = 47694
$ws.Range("J111").Value = 47694
$ws.Range("L111").Value = 47694
$ws.Range("N111").Value = -55874
$ws.Range("H117").Value = 49914
$ws.Range("J117").Value = 49914
$ws.Range("L117").Value = 49914
$ws.Range("N117").Value = -59092
$ws.Range("H130").Value = 49077.75
$ws.Range("J130").Value = 49077.75
$ws.Range("L130").Value = 49077.75
$ws.Range("N130").Value = -59117.75
$ws.Range("H132").Value = 50336.152
$ws.Range("J132").Value = 50336.152
$ws.Range("L132").Value = 50336.152
$ws.Range("N132").Value = -60456.152

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49780
$ws.Range("J20").Value = 49780
$ws.Range("L20").Value = 49780
$ws.Range("N20").Value = -50252
$ws.Range("H30").Value = 49780
$ws.Range("J30").Value = 49780
$ws.Range("L30").Value = 49780
$ws.Range("N30").Value = -49962
$ws.Range("H31").Value = 4677.7075
$ws.Range("I31").Value = 1575
$ws.Range("J31").Value = 6616.9
$ws.Range("K31").Value = 1575
$ws.Range("L31").Value = 6616.9
$ws.Range("M31").Value = -1280
$ws.Range("N31").Value = -7206.9
$ws.Range("H34").Value = 4677.7075
$ws.Range("I34").Value = 1575
$ws.Range("J34").Value = 6616.9
$ws.Range("K34").Value = 1575
$ws.Range("L34").Value = 6616.9
$ws.Range("M34").Value = -1373
$ws.Range("N34").Value = -7020.9
$ws.Range("H100").Value = 46850.668
$ws.Range("J100").Value = 46850.668
$ws.Range("L100").Value = 46850.668
$ws.Range("N100").Value = -49014.668
$ws.Range("H110").Value = 39697.5
$ws.Range("J110").Value = 39697.5
$ws.Range("L110").Value = 39697.5
$ws.Range("N110").Value = -47877.5
$ws.Range("H111").Value = 46988
$ws.Range("J111").Value = 46988
$ws.Range("L111").Value = 46988
$ws.Range("N111").Value = -55168
$ws.Range("H112").Value = 37360.43
$ws.Range("J112").Value = 37360.43
$ws.Range("L112").Value = 37360.43
$ws.Range("N112").Value = -40314.43
$ws.Range("H116").Value = 47822.332
$ws.Range("J116").Value = 47822.332
$ws.Range("L116").Value = 47822.332
$ws.Range("N116").Value = -57000.332
$ws.Range("H119").Value = 48581.668
$ws.Range("J119").Value = 48581.668
$ws.Range("L119").Value = 48581.668
$ws.Range("N119").Value = -58257.668
$ws.Range("H128").Value = 49780
$ws.Range("J128").Value = 49780
$ws.Range("L128").Value = 49780
$ws.Range("N128").Value = -59740
$ws.Range("H135").Value = 38694.668
$ws.Range("J135").Value = 38694.668
$ws.Range("L135").Value = 38694.668
$ws.Range("N135").Value = -48834.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 12500155
$ws.Range("I2").Value = 51.166668
$ws.Range("J2").Value = 31250310
$ws.Range("K2").Value = 307.000008
$ws.Range("L2").Value = 187501860
$ws.Range("M2").Value = -194.000008
$ws.Range("N2").Value = -187502086
$ws.Range("H34").Value = 1342.9131
$ws.Range("J34").Value = 1394.8636
$ws.Range("L34").Value = 4184.5908
$ws.Range("N34").Value = -4352.5908
$ws.Range("H39").Value = 3000
$ws.Range("J39").Value = 3666.6667
$ws.Range("L39").Value = 11000.0001
$ws.Range("N39").Value = -11588.0001
$ws.Range("H109").Value = 13279.556
$ws.Range("I109").Value = 41818.285
$ws.Range("J109").Value = 3291
$ws.Range("K109").Value = 125454.855
$ws.Range("L109").Value = 9873
$ws.Range("M109").Value = -124414.855
$ws.Range("N109").Value = -11953
$ws.Range("H138").Value = 2635.65
$ws.Range("I138").Value = 2040.4375
$ws.Range("J138").Value = 5016.5
$ws.Range("K138").Value = 6121.3125
$ws.Range("L138").Value = 15049.5
$ws.Range("M138").Value = -981.3125
$ws.Range("N138").Value = -25329.5
$ws.Range("H140").Value = 69995.734
$ws.Range("I140").Value = 102599.7
$ws.Range("J140").Value = 4787.8
$ws.Range("K140").Value = 307799.1
$ws.Range("L140").Value = 14363.4
$ws.Range("M140").Value = -302619.1
$ws.Range("N140").Value = -24723.4
$ws.Range("H141").Value = 33336768
$ws.Range("I141").Value = 43481324
$ws.Range("J141").Value = 4655.5713
$ws.Range("K141").Value = 130443972
$ws.Range("L141").Value = 13966.7139
$ws.Range("M141").Value = -130438792
$ws.Range("N141").Value = -24326.7139

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 90009
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H27").Value = 41000
$ws.Range("J27").Value = 41000
$ws.Range("L27").Value = 41000
$ws.Range("N27").Value = -41332
$ws.Range("H31").Value = 39000
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 39000
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 39000
$ws.Range("N31").Value = -39584
$ws.Range("M31").ClearContents()
$ws.Range("H37").Value = 39000
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 39000
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 39000
$ws.Range("N37").Value = -39554
$ws.Range("M37").ClearContents()
$ws.Range("H64").Value = 29833.334
$ws.Range("J64").Value = 29833.334
$ws.Range("L64").Value = 29833.334
$ws.Range("N64").Value = -30329.334
$ws.Range("H67").Value = 29833.334
$ws.Range("J67").Value = 29833.334
$ws.Range("L67").Value = 29833.334
$ws.Range("N67").Value = -31549.334
$ws.Range("H110").Value = 39351
$ws.Range("J110").Value = 39351
$ws.Range("L110").Value = 39351
$ws.Range("N110").Value = -47531
$ws.Range("H114").Value = 48712
$ws.Range("J114").Value = 48712
$ws.Range("L114").Value = 48712
$ws.Range("N114").Value = -57390
$ws.Range("H119").Value = 48753
$ws.Range("J119").Value = 48753
$ws.Range("L119").Value = 48753
$ws.Range("N119").Value = -58429
$ws.Range("H126").Value = 6260.0835
$ws.Range("I126").Value = 11172
$ws.Range("J126").Value = 2103.8462
$ws.Range("K126").Value = 33516
$ws.Range("L126").Value = 6311.5386
$ws.Range("M126").Value = -31046
$ws.Range("N126").Value = -11251.5386
$ws.Range("H130").Value = 52986
$ws.Range("J130").Value = 52986
$ws.Range("L130").Value = 52986
$ws.Range("N130").Value = -63026
$ws.Range("H134").Value = 24921.428
$ws.Range("J134").Value = 24921.428
$ws.Range("L134").Value = 74764.284
$ws.Range("N134").Value = -79834.284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2895.111
$ws.Range("I7").Value = 2176
$ws.Range("J7").Value = 4333.3335
$ws.Range("K7").Value = 2176
$ws.Range("L7").Value = 4333.3335
$ws.Range("M7").Value = -2064
$ws.Range("N7").Value = -4557.3335
$ws.Range("H36").Value = 46707
$ws.Range("J36").Value = 46707
$ws.Range("L36").Value = 46707
$ws.Range("N36").Value = -47831
$ws.Range("H102").Value = 48545
$ws.Range("J102").Value = 48545
$ws.Range("L102").Value = 48545
$ws.Range("N102").Value = -55035
$ws.Range("H108").Value = 41081.332
$ws.Range("J108").Value = 41081.332
$ws.Range("L108").Value = 41081.332
$ws.Range("N108").Value = -48761.332
$ws.Range("H111").Value = 35954.332
$ws.Range("J111").Value = 35954.332
$ws.Range("L111").Value = 35954.332
$ws.Range("N111").Value = -44134.332
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H116").Value = 50453.332
$ws.Range("J116").Value = 50453.332
$ws.Range("L116").Value = 50453.332
$ws.Range("N116").Value = -59631.332
$ws.Range("H117").Value = 43594.668
$ws.Range("J117").Value = 43594.668
$ws.Range("L117").Value = 43594.668
$ws.Range("N117").Value = -52772.668
$ws.Range("H118").Value = 43401
$ws.Range("J118").Value = 43401
$ws.Range("L118").Value = 43401
$ws.Range("N118").Value = -46715
$ws.Range("H119").Value = 45673.332
$ws.Range("J119").Value = 45673.332
$ws.Range("L119").Value = 45673.332
$ws.Range("N119").Value = -55349.332
$ws.Range("H120").Value = 56015
$ws.Range("J120").Value = 56015
$ws.Range("L120").Value = 56015
$ws.Range("N120").Value = -65691
$ws.Range("H121").Value = 37940.668
$ws.Range("J121").Value = 37940.668
$ws.Range("L121").Value = 37940.668
$ws.Range("N121").Value = -41434.668
$ws.Range("H124").Value = 39498
$ws.Range("J124").Value = 39498
$ws.Range("L124").Value = 39498
$ws.Range("N124").Value = -49318
$ws.Range("H125").Value = 49707
$ws.Range("J125").Value = 49707
$ws.Range("L125").Value = 49707
$ws.Range("N125").Value = -59547
$ws.Range("H126").Value = 2895.111
$ws.Range("I126").Value = 2176
$ws.Range("J126").Value = 4333.3335
$ws.Range("K126").Value = 6528
$ws.Range("L126").Value = 13000.0005
$ws.Range("M126").Value = -4058
$ws.Range("N126").Value = -17940.0005
$ws.Range("H127").Value = 48190
$ws.Range("J127").Value = 48190
$ws.Range("L127").Value = 48190
$ws.Range("N127").Value = -58110
$ws.Range("H128").Value = 44714.5
$ws.Range("J128").Value = 44714.5
$ws.Range("L128").Value = 44714.5
$ws.Range("N128").Value = -54674.5
$ws.Range("H130").Value = 48407
$ws.Range("J130").Value = 48407
$ws.Range("L130").Value = 48407
$ws.Range("N130").Value = -58447

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 35813.332
$ws.Range("J92").Value = 35813.332
$ws.Range("L92").Value = 35813.332
$ws.Range("N92").Value = -40805.332
$ws.Range("H93").Value = 37315
$ws.Range("J93").Value = 37315
$ws.Range("L93").Value = 37315
$ws.Range("N93").Value = -42307
$ws.Range("H99").Value = 39650.668
$ws.Range("J99").Value = 41998.855
$ws.Range("L99").Value = 41998.855
$ws.Range("N99").Value = -47988.855
$ws.Range("H102").Value = 29582.5
$ws.Range("J102").Value = 29582.5
$ws.Range("L102").Value = 29582.5
$ws.Range("N102").Value = -36072.5
$ws.Range("H108").Value = 22818
$ws.Range("J108").Value = 22818
$ws.Range("L108").Value = 22818
$ws.Range("N108").Value = -30498
$ws.Range("H109").Value = 33804.5
$ws.Range("J109").Value = 33804.5
$ws.Range("L109").Value = 33804.5
$ws.Range("N109").Value = -36578.5
$ws.Range("H110").Value = 49644
$ws.Range("J110").Value = 49644
$ws.Range("L110").Value = 49644
$ws.Range("N110").Value = -57824
$ws.Range("H117").Value = 42774.75
$ws.Range("J117").Value = 42774.75
$ws.Range("L117").Value = 42774.75
$ws.Range("N117").Value = -51952.75
$ws.Range("H120").Value = 42203.668
$ws.Range("J120").Value = 42203.668
$ws.Range("L120").Value = 42203.668
$ws.Range("N120").Value = -51879.668
$ws.Range("H123").Value = 41837.75
$ws.Range("J123").Value = 41837.75
$ws.Range("L123").Value = 41837.75
$ws.Range("N123").Value = -51637.75
$ws.Range("H131").Value = 48711
$ws.Range("J131").Value = 48711
$ws.Range("L131").Value = 48711
$ws.Range("N131").Value = -58791
$ws.Range("H133").Value = 83557
$ws.Range("J133").Value = 83557
$ws.Range("L133").Value = 83557
$ws.Range("N133").Value = -93677
